$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F29").Value = 23206
$ws.Range("F30").Value = 12048
$ws.Range("F31").Value = 13302
$ws.Range("B32").Value = 0.8092818759159746
$ws.Range("G32").Value = 103078826.9229743
$ws.Range("M32").Value = 1331385835.325045
$ws.Range("F35").Value = 24991
